# Auto-generated edit script reproducing the commit:
# "Atualizacao de bases das ligas, do dia: 28-04-2024 as 23:19"
#
# Summary of the edit:
#  - The shared strings "Lion City Sailors FC" and "DPMM FC" swap meaning
#    (a data-correction: these two teams were mislabeled in the source feed).
#  - Six pairs of data rows (6/7, 8/9, 22/23, 26/27, 47/48, 54/55) had their
#    entire row content exchanged (every column except "A", the running index).
#  - A handful of standalone HomeTeam/AwayTeam cells elsewhere in the sheet
#    also flip between these two team names as a consequence of the swap.
#
# Rather than re-deriving this at runtime, every affected cell is set directly
# to its final, known value (computed from the unified diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 ---
$ws.Range("F2").Value = 'DPMM FC'

# --- row 6 ---
$ws.Range("B6").Value = 6228587
$ws.Range("F6").Value = 'Hougang United FC'
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 'A'
$ws.Range("K6").Value = 3.75
$ws.Range("L6").Value = 2.3
$ws.Range("M6").Value = 2.4
$ws.Range("N6").Value = 3.75
$ws.Range("O6").Value = 2.3
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 1.975
$ws.Range("R6").Value = 1.875
$ws.Range("S6").Value = 3.5
$ws.Range("T6").Value = 1.85
$ws.Range("U6").Value = 2
$ws.Range("W6").Value = -1
$ws.Range("X6").Value = 1.3
$ws.Range("Y6").Value = -1
$ws.Range("Z6").Value = 0.875
$ws.Range("AA6").Value = 0.8500000000000001
$ws.Range("AB6").Value = -1

# --- row 7 ---
$ws.Range("B7").Value = 6228027
$ws.Range("F7").Value = 'Tampines Rovers FC'
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 'D'
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = 2.25
$ws.Range("M7").Value = 2.7
$ws.Range("N7").Value = 4
$ws.Range("O7").Value = 2.05
$ws.Range("P7").Value = 0.25
$ws.Range("Q7").Value = 1.925
$ws.Range("R7").Value = 1.925
$ws.Range("S7").Value = 4
$ws.Range("T7").Value = 1.95
$ws.Range("U7").Value = 1.9
$ws.Range("W7").Value = 3
$ws.Range("X7").Value = -1
$ws.Range("Y7").Value = 0.4625
$ws.Range("Z7").Value = -0.5
$ws.Range("AA7").Value = -1
$ws.Range("AB7").Value = 0.8999999999999999

# --- row 8 ---
$ws.Range("B8").Value = 6228028
$ws.Range("E8").Value = 'DPMM FC'
$ws.Range("F8").Value = 'Young Lions'
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 1
$ws.Range("J8").Value = 1.083
$ws.Range("K8").Value = 9
$ws.Range("L8").Value = 15
$ws.Range("M8").Value = 1.142
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 12
$ws.Range("P8").Value = -2.25
$ws.Range("Q8").Value = 1.85
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 4
$ws.Range("T8").Value = 1.875
$ws.Range("U8").Value = 1.975
$ws.Range("V8").Value = 0.1419999999999999
$ws.Range("Y8").Value = 0.8500000000000001
$ws.Range("AA8").Value = 0.875
$ws.Range("AB8").Value = -1

# --- row 9 ---
$ws.Range("B9").Value = 6228588
$ws.Range("E9").Value = 'Tampines Rovers FC'
$ws.Range("F9").Value = 'Hougang United FC'
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 1.333
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = 6.25
$ws.Range("M9").Value = 1.615
$ws.Range("N9").Value = 4.5
$ws.Range("O9").Value = 3.8
$ws.Range("P9").Value = -1
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 1.8
$ws.Range("S9").Value = 3.5
$ws.Range("T9").Value = 1.85
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 0.615
$ws.Range("Y9").Value = 1.05
$ws.Range("AA9").Value = -1
$ws.Range("AB9").Value = 1

# --- row 11 ---
$ws.Range("E11").Value = 'Lion City Sailors FC'

# --- row 13 ---
$ws.Range("E13").Value = 'Lion City Sailors FC'

# --- row 15 ---
$ws.Range("E15").Value = 'DPMM FC'

# --- row 17 ---
$ws.Range("F17").Value = 'Lion City Sailors FC'

# --- row 18 ---
$ws.Range("F18").Value = 'DPMM FC'

# --- row 22 ---
$ws.Range("B22").Value = 6228600
$ws.Range("E22").Value = 'Tanjong Pagar United'
$ws.Range("F22").Value = 'Balestier Khalsa FC'
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 3
$ws.Range("J22").Value = 3.2
$ws.Range("K22").Value = 4
$ws.Range("L22").Value = 1.8
$ws.Range("M22").Value = 3.4
$ws.Range("N22").Value = 4.2
$ws.Range("O22").Value = 1.8
$ws.Range("P22").Value = 0.75
$ws.Range("Q22").Value = 1.825
$ws.Range("R22").Value = 2.025
$ws.Range("S22").Value = 4.5
$ws.Range("X22").Value = 0.8
$ws.Range("Y22").Value = -0.5
$ws.Range("Z22").Value = 0.5125
$ws.Range("AA22").Value = 1
$ws.Range("AB22").Value = -1

# --- row 23 ---
$ws.Range("B23").Value = 6228599
$ws.Range("E23").Value = 'Hougang United FC'
$ws.Range("F23").Value = 'Tampines Rovers FC'
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1
$ws.Range("J23").Value = 4.75
$ws.Range("K23").Value = 4.2
$ws.Range("L23").Value = 1.5
$ws.Range("M23").Value = 7.5
$ws.Range("N23").Value = 4.75
$ws.Range("O23").Value = 1.3
$ws.Range("P23").Value = 1.5
$ws.Range("Q23").Value = 2
$ws.Range("R23").Value = 1.85
$ws.Range("S23").Value = 3.75
$ws.Range("X23").Value = 0.3
$ws.Range("Y23").Value = 1
$ws.Range("Z23").Value = -1
$ws.Range("AA23").Value = -1
$ws.Range("AB23").Value = 0.8500000000000001

# --- row 24 ---
$ws.Range("F24").Value = 'DPMM FC'

# --- row 25 ---
$ws.Range("E25").Value = 'Lion City Sailors FC'

# --- row 26 ---
$ws.Range("B26").Value = 6228602
$ws.Range("E26").Value = 'Tampines Rovers FC'
$ws.Range("F26").Value = 'Tanjong Pagar United'
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 'H'
$ws.Range("M26").Value = 1.125
$ws.Range("O26").Value = 17
$ws.Range("P26").Value = -2.5
$ws.Range("Q26").Value = 1.9
$ws.Range("R26").Value = 1.95
$ws.Range("S26").Value = 4.25
$ws.Range("T26").Value = 1.975
$ws.Range("U26").Value = 1.875
$ws.Range("V26").Value = 0.125
$ws.Range("X26").Value = -1
$ws.Range("Z26").Value = 0.95
$ws.Range("AA26").Value = -1
$ws.Range("AB26").Value = 0.875

# --- row 27 ---
$ws.Range("B27").Value = 6228032
$ws.Range("E27").Value = 'DPMM FC'
$ws.Range("F27").Value = 'Lion City Sailors FC'
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 'A'
$ws.Range("M27").Value = 1.142
$ws.Range("O27").Value = 13
$ws.Range("P27").Value = -2.25
$ws.Range("Q27").Value = 1.825
$ws.Range("R27").Value = 2.025
$ws.Range("S27").Value = 4
$ws.Range("T27").Value = 1.825
$ws.Range("U27").Value = 2.025
$ws.Range("V27").Value = -1
$ws.Range("X27").Value = 12
$ws.Range("Z27").Value = 1.025
$ws.Range("AA27").Value = 0
$ws.Range("AB27").Value = 0

# --- row 32 ---
$ws.Range("E32").Value = 'Lion City Sailors FC'

# --- row 33 ---
$ws.Range("F33").Value = 'DPMM FC'

# --- row 34 ---
$ws.Range("E34").Value = 'DPMM FC'

# --- row 38 ---
$ws.Range("E38").Value = 'Lion City Sailors FC'

# --- row 41 ---
$ws.Range("E41").Value = 'DPMM FC'

# --- row 43 ---
$ws.Range("F43").Value = 'DPMM FC'

# --- row 44 ---
$ws.Range("F44").Value = 'Lion City Sailors FC'

# --- row 47 ---
$ws.Range("B47").Value = 6228618
$ws.Range("E47").Value = 'Young Lions'
$ws.Range("F47").Value = 'Lion City Sailors FC'
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 2
$ws.Range("J47").Value = 5
$ws.Range("K47").Value = 5
$ws.Range("L47").Value = 1.4
$ws.Range("M47").Value = 5
$ws.Range("N47").Value = 5.25
$ws.Range("O47").Value = 1.4
$ws.Range("P47").Value = 1.25
$ws.Range("Q47").Value = 2
$ws.Range("R47").Value = 1.85
$ws.Range("S47").Value = 4
$ws.Range("T47").Value = 1.975
$ws.Range("U47").Value = 1.875
$ws.Range("X47").Value = 0.3999999999999999
$ws.Range("Z47").Value = 0.8500000000000001
$ws.Range("AA47").Value = -1
$ws.Range("AB47").Value = 0.875

# --- row 48 ---
$ws.Range("B48").Value = 6228619
$ws.Range("E48").Value = 'Geylang International'
$ws.Range("F48").Value = 'Balestier Khalsa FC'
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 6
$ws.Range("J48").Value = 2.1
$ws.Range("K48").Value = 4
$ws.Range("L48").Value = 2.55
$ws.Range("M48").Value = 2.15
$ws.Range("N48").Value = 3.6
$ws.Range("O48").Value = 2.8
$ws.Range("P48").Value = -0.25
$ws.Range("Q48").Value = 1.95
$ws.Range("R48").Value = 1.9
$ws.Range("S48").Value = 4.5
$ws.Range("T48").Value = 2
$ws.Range("U48").Value = 1.85
$ws.Range("X48").Value = 1.8
$ws.Range("Z48").Value = 0.8999999999999999
$ws.Range("AA48").Value = 1
$ws.Range("AB48").Value = -1

# --- row 49 ---
$ws.Range("F49").Value = 'DPMM FC'

# --- row 50 ---
$ws.Range("F50").Value = 'Lion City Sailors FC'

# --- row 54 ---
$ws.Range("B54").Value = 7094656
$ws.Range("E54").Value = 'Tanjong Pagar United'
$ws.Range("F54").Value = 'Lion City Sailors FC'
$ws.Range("H54").Value = 1
$ws.Range("I54").Value = 'D'
$ws.Range("J54").Value = 2.15
$ws.Range("K54").Value = 3.75
$ws.Range("L54").Value = 2.7
$ws.Range("M54").Value = 2.1
$ws.Range("N54").Value = 4.2
$ws.Range("O54").Value = 2.625
$ws.Range("P54").Value = -0.25
$ws.Range("Q54").Value = 1.925
$ws.Range("R54").Value = 1.925
$ws.Range("S54").Value = 4.25
$ws.Range("T54").Value = 1.9
$ws.Range("U54").Value = 1.95
$ws.Range("W54").Value = 3.2
$ws.Range("X54").Value = -1
$ws.Range("Y54").Value = -0.5
$ws.Range("Z54").Value = 0.4625
$ws.Range("AB54").Value = 0.95

# --- row 55 ---
$ws.Range("B55").Value = 7098763
$ws.Range("E55").Value = 'Balestier Khalsa FC'
$ws.Range("F55").Value = 'Tampines Rovers FC'
$ws.Range("H55").Value = 3
$ws.Range("I55").Value = 'A'
$ws.Range("J55").Value = 5.25
$ws.Range("K55").Value = 4.2
$ws.Range("L55").Value = 1.5
$ws.Range("M55").Value = 5
$ws.Range("N55").Value = 4.5
$ws.Range("O55").Value = 1.45
$ws.Range("P55").Value = 1.25
$ws.Range("Q55").Value = 2
$ws.Range("R55").Value = 1.85
$ws.Range("S55").Value = 5
$ws.Range("T55").Value = 1.925
$ws.Range("U55").Value = 1.925
$ws.Range("W55").Value = -1
$ws.Range("X55").Value = 0.45
$ws.Range("Y55").Value = -1
$ws.Range("Z55").Value = 0.8500000000000001
$ws.Range("AB55").Value = 0.925

# --- row 57 ---
$ws.Range("E57").Value = 'DPMM FC'

